$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 5 (2013年) values ---
$ws.Range("B5").Value = 8517893
$ws.Range("C5").Value = 626211
$ws.Range("D5").Value = 226528
$ws.Range("E5").Value = 86359
$ws.Range("F5").Value = 84466
$ws.Range("G5").Value = 7271891
$ws.Range("H5").Value = 222438

# --- Append new row 12 (2021年) ---
$ws.Range("A12").Value = "2021年"
$ws.Range("B12").Value = 28665212
$ws.Range("C12").Value = 365021
$ws.Range("D12").Value = 323277
$ws.Range("E12").Value = 109043
$ws.Range("F12").Value = 135201
$ws.Range("G12").Value = 27545002
$ws.Range("H12").Value = 187668

# Match the formatting used by the other "year" cells in column A
# (bold, centered, thin-bordered style) by copying formats from A11,
# which keeps the same style index instead of creating a new one.
$ws.Range("A11").Copy()
$ws.Range("A12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Re-assert the value/type (PasteSpecial formats-only shouldn't touch it,
# but make sure the text is exactly what we want).
$ws.Range("A12").Value = "2021年"
